# Update "想去人数" (F column) figures to the newly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 385
$ws.Range("F4").Value = 433
$ws.Range("F5").Value = 1174
$ws.Range("F8").Value = 1178
$ws.Range("F9").Value = 1665
$ws.Range("F10").Value = 6179
$ws.Range("F12").Value = 1792
$ws.Range("F13").Value = 466
$ws.Range("F15").Value = 7
$ws.Range("F16").Value = 2
$ws.Range("F18").Value = 8
$ws.Range("F19").Value = 6471
$ws.Range("F21").Value = 51
$ws.Range("F22").Value = 162
$ws.Range("F23").Value = 102
$ws.Range("F24").Value = 1694
$ws.Range("F26").Value = 10
$ws.Range("F28").Value = 156
$ws.Range("F29").Value = 1518
$ws.Range("F30").Value = 742
$ws.Range("F31").Value = 296
$ws.Range("F36").Value = 3885

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 328
$ws.Range("F5").Value = 188
$ws.Range("F8").Value = 419
$ws.Range("F19").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9526
$ws.Range("F3").Value = 2249
$ws.Range("F5").Value = 227

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9526
$ws.Range("F3").Value = 2249
$ws.Range("F5").Value = 385
$ws.Range("F6").Value = 433
$ws.Range("F7").Value = 1174
$ws.Range("F11").Value = 328
$ws.Range("F12").Value = 1178
$ws.Range("F13").Value = 227
$ws.Range("F14").Value = 1665
$ws.Range("F15").Value = 6179
$ws.Range("F17").Value = 1792
$ws.Range("F19").Value = 466
$ws.Range("F21").Value = 7
$ws.Range("F23").Value = 6471
$ws.Range("F25").Value = 51
$ws.Range("F26").Value = 162
$ws.Range("F27").Value = 102
$ws.Range("F28").Value = 1694
$ws.Range("F31").Value = 156
$ws.Range("F32").Value = 1518
$ws.Range("F33").Value = 742
$ws.Range("F35").Value = 296
$ws.Range("F43").Value = 4
$ws.Range("F44").Value = 3885

